# Scheduled-runner style refresh of market/profit figures across the
# Marilith_Profits workbook's per-class leve sheets (ALC, ARM, BSM, CRP,
# CUL, GSM, LTW, WVR). Updates currentAveragePrice* / LevePrice* /
# LeveProfit* columns (H:N) with freshly pulled values; some rows no
# longer have a computed profit cell and those are cleared outright.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 704.8946999999999
$ws.Range("I28").Value = 733.7059
$ws.Range("K28").Value = 733.7059
$ws.Range("M28").Value = -248.7059

$ws.Range("H93").Value = 67500
$ws.Range("J93").Value = 67500
$ws.Range("L93").Value = 67500
$ws.Range("N93").Value = -72492

$ws.Range("H125").Value = 4899.25
$ws.Range("I125").Value = 2532.3333
$ws.Range("K125").Value = 22790.9997
$ws.Range("M125").Value = -20330.9997

$ws.Range("H141").Value = 3616.0833
$ws.Range("I141").Value = 4033.5715
$ws.Range("J141").Value = 693.6667
$ws.Range("K141").Value = 12100.7145
$ws.Range("L141").Value = 2081.0001
$ws.Range("M141").Value = -6920.7145
$ws.Range("N141").Value = -12441.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()

$ws.Range("H61").Value = 4160
$ws.Range("I61").Value = 2742.6
$ws.Range("K61").Value = 2742.6
$ws.Range("M61").Value = -2530.6

$ws.Range("H94").Value = 84999.5
$ws.Range("J94").Value = 84999.5
$ws.Range("L94").Value = 84999.5
$ws.Range("N94").Value = -86801.5

$ws.Range("H132").Value = 3038.76
$ws.Range("J132").Value = 9398.4
$ws.Range("L132").Value = 28195.2
$ws.Range("N132").Value = -33255.2

$ws.Range("H134").Value = 50000
$ws.Range("J134").Value = 50000
$ws.Range("L134").Value = 50000
$ws.Range("N134").Value = -60140

$ws.Range("H136").Value = 4160
$ws.Range("I136").Value = 2742.6
$ws.Range("K136").Value = 8227.799999999999
$ws.Range("M136").Value = -5677.799999999999

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 19876
$ws.Range("J76").Value = 19876
$ws.Range("L76").Value = 19876
$ws.Range("N76").Value = -20506

$ws.Range("H79").Value = 19876
$ws.Range("J79").Value = 19876
$ws.Range("L79").Value = 19876
$ws.Range("N79").Value = -22060

$ws.Range("H88").Value = 15000
$ws.Range("I88").Value = 10000
$ws.Range("J88").Value = 20000
$ws.Range("K88").Value = 10000
$ws.Range("L88").Value = 20000
$ws.Range("M88").Value = -9594
$ws.Range("N88").Value = -20812

$ws.Range("H91").Value = 15000
$ws.Range("I91").Value = 10000
$ws.Range("J91").Value = 20000
$ws.Range("K91").Value = 10000
$ws.Range("L91").Value = 20000
$ws.Range("M91").Value = -8596
$ws.Range("N91").Value = -22808

$ws.Range("H92").Value = 50200
$ws.Range("J92").Value = 50200
$ws.Range("L92").Value = 50200
$ws.Range("N92").Value = -55192

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4376
$ws.Range("I31").Value = 4174.5
$ws.Range("K31").Value = 4174.5
$ws.Range("M31").Value = -3879.5

$ws.Range("H34").Value = 4376
$ws.Range("I34").Value = 4174.5
$ws.Range("K34").Value = 4174.5
$ws.Range("M34").Value = -3972.5

$ws.Range("H92").Value = 55999.4
$ws.Range("J92").Value = 55999.4
$ws.Range("L92").Value = 55999.4
$ws.Range("N92").Value = -60991.4

$ws.Range("H132").Value = 1781.6666
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

$ws.Range("H75").Value = 654.2
$ws.Range("I75").Value = 692.5
$ws.Range("J75").Value = 628.6667
$ws.Range("K75").Value = 2077.5
$ws.Range("L75").Value = 1886.0001
$ws.Range("M75").Value = -1079.5
$ws.Range("N75").Value = -3882.0001

$ws.Range("H78").Value = 654.2
$ws.Range("I78").Value = 692.5
$ws.Range("J78").Value = 628.6667
$ws.Range("K78").Value = 6232.5
$ws.Range("L78").Value = 5658.0003
$ws.Range("M78").Value = -1240.5
$ws.Range("N78").Value = -15642.0003

$ws.Range("H123").Value = 999.5
$ws.Range("J123").Value = 999
$ws.Range("L123").Value = 2997
$ws.Range("N123").Value = -7897

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 5625
$ws.Range("J92").Value = 5625
$ws.Range("L92").Value = 5625
$ws.Range("N92").Value = -9369

$ws.Range("H93").Value = 38000
$ws.Range("J93").Value = 38000
$ws.Range("L93").Value = 38000
$ws.Range("N93").Value = -41744

$ws.Range("H133").Value = 100707.75
$ws.Range("J133").Value = 100707.75
$ws.Range("L133").Value = 100707.75
$ws.Range("N133").Value = -110827.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 884.2
$ws.Range("I22").Value = 904.6667
$ws.Range("K22").Value = 904.6667
$ws.Range("M22").Value = -609.6667

$ws.Range("H27").Value = 884.2
$ws.Range("I27").Value = 904.6667
$ws.Range("K27").Value = 904.6667
$ws.Range("M27").Value = -797.6667

$ws.Range("H31").Value = 1345
$ws.Range("I31").Value = 1183
$ws.Range("K31").Value = 1183
$ws.Range("M31").Value = -935

$ws.Range("H55").Value = 885.7857
$ws.Range("J55").Value = 669
$ws.Range("L55").Value = 669
$ws.Range("N55").Value = -1015

$ws.Range("H134").Value = 60000
$ws.Range("J134").Value = 60000
$ws.Range("L134").Value = 60000
$ws.Range("N134").Value = -70140

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3053.2
$ws.Range("I126").Value = 3002.7693
$ws.Range("K126").Value = 9008.3079
$ws.Range("M126").Value = -6538.3079

$ws.Range("H136").Value = 2919.5789
$ws.Range("I136").Value = 2804
$ws.Range("K136").Value = 8412
$ws.Range("M136").Value = -5862

$ws.Range("H137").Value = 250000
$ws.Range("J137").Value = 250000
$ws.Range("L137").Value = 250000
$ws.Range("N137").Value = -260200

Write-Output "Marilith_Profits sheets refreshed"
